$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197, pushing existing rows 197-272 down to 198-273.
$ws.Rows("197:197").Insert()

# Populate the new row 197 with the new record (copy of neighboring row structure
# with updated Fecha/Volumen/Precio promedio ponderado/Precio $/Kg values).
$ws.Cells.Item(197, 1).Value = 10
$ws.Cells.Item(197, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(197, 3).Value = "La Araucanía"
$ws.Cells.Item(197, 4).Value = 44784
$ws.Cells.Item(197, 5).Value = 9
$ws.Cells.Item(197, 6).Value = 100112039
$ws.Cells.Item(197, 7).Value = "Ciboulette"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 80
$ws.Cells.Item(197, 11).Value = 5000
$ws.Cells.Item(197, 12).Value = 6000
$ws.Cells.Item(197, 13).Value = 5500
$ws.Cells.Item(197, 14).Value = "`$/docena de atados"
$ws.Cells.Item(197, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(197, 16).Value = 1833
$ws.Cells.Item(197, 17).Value = 3
$ws.Cells.Item(197, 18).Value = "Hortaliza"
